# Updated header files so that the app shows a different header depending on
# whether you're logged in; added login/register text to main page; updated
# testing spreadsheet.
#
# This script applies the "testing spreadsheet" portion of that change to
# the moodTube Testing Spreadsheet workbook (Remediation table on Sheet1):
# four new bugs were logged (rows 63-66), two existing rows picked up an
# "Assigned" owner (and a completed date for row 60), and the stale Test ID
# on row 61 was cleared out.
#
# Note: the workbook was also re-saved from a "v1" folder into a "v2"
# folder on the author's machine, which is why the diff touches the
# Microsoft-internal <x15ac:absPath> breadcrumb in xl/workbook.xml and the
# sheetView's topLeftCell scroll position. Neither of those is a real
# document property - Excel writes them for its own bookkeeping and does
# not expose them anywhere in the Application/Workbook/Window object model
# (no VBA/COM property maps to either one), so they can't be set from
# automation code here; only the cell data/view-selection changes below
# are reproducible that way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 60: add Assigned / Completed Date -------------------------------
$ws.Range("E60").Value = "Sarah"
# Copy the date number format used by the other "Completed Date" cells
# (e.g. F59) so the new cell gets the same style, then set its value.
$ws.Range("F59").Copy()
$ws.Range("F60").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F60").Value2 = 43071

# --- Row 61: the "Test ID" column (B61, "?") is no longer applicable -----
$ws.Range("B61").ClearContents()

# --- Row 62: add Assigned -------------------------------------------------
$ws.Range("E62").Value = "Sarah"

# --- New row 63 -------------------------------------------------------
$ws.Range("A63").Value2 = 12
$ws.Range("C63").Value = "system allows user to enter any value for email when that's not correct"
$ws.Range("D63").Value2 = 3
$ws.Range("E63").Value = "Sarah"

# --- New row 64 -------------------------------------------------------
$ws.Range("A64").Value2 = 13
$ws.Range("C64").Value = "clicking like dislike share or unlike or undislike reroutes the user instead of letting them continue the search"
$ws.Range("D64").Value2 = 2

# --- New row 65 -------------------------------------------------------
$ws.Range("A65").Value2 = 14
$ws.Range("C65").Value = "buttons for logged-in only functions show for unlogged in users"
$ws.Range("D65").Value2 = 2
$ws.Range("E65").Value = "Sonya"
$ws.Range("F59").Copy()
$ws.Range("F65").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F65").Value2 = 43071

# --- New row 66 -------------------------------------------------------
$ws.Range("A66").Value2 = 15
$ws.Range("C66").Value = "share button doesn't work"
$ws.Range("D66").Value2 = 2
$ws.Range("E66").Value = "Jordan"

$excel.CutCopyMode = $false

# --- Update the view so the selection matches the saved workbook ---------
$ws.Range("C62").Select()
$win = $wb.Windows.Item(1)
$win.ScrollRow = 46
$win.ScrollColumn = 1
